$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A71").Value = "2025-04-29 10:36:47"
$ws.Range("B71").Value = 214
